$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change existing value in I13 from "SHEF" to "sfawfh"
$ws.Range("I13").Value = "sfawfh"

# Add new values
$ws.Range("M12").Value = "oAHSFIyfe"
$ws.Range("K16").Value = "iuhefiuHWEF;"

# Update selection to match the final edited cell
$ws.Range("K16").Select()
